$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before row 145 (shifts 145:237 -> 147:239,
# and duplicates the old bottom two rows into the new bottom two rows 238:239
# automatically since they are simply pushed down).
$ws.Rows("145:146").Insert()

# Fill in the two newly-inserted rows (145 and 146) with their data. All of
# the "static" columns (A, B, C, E, F, G, H, N, Q, R) are identical to every
# other data row in this block, so copy them from row 147 (the row that used
# to be row 145 before the insert) to stay consistent, then set the columns
# that actually carry new data for this entry.

$ws.Cells.Item(145, 1).Value  = $ws.Cells.Item(147, 1).Value()   # A - Mercado ID
$ws.Cells.Item(145, 2).Value  = $ws.Cells.Item(147, 2).Value()   # B - Mercado
$ws.Cells.Item(145, 3).Value  = $ws.Cells.Item(147, 3).Value()   # C - Region
$ws.Cells.Item(145, 4).Value  = 44452                             # D - Fecha
$ws.Cells.Item(145, 5).Value  = $ws.Cells.Item(147, 5).Value()   # E - Codreg
$ws.Cells.Item(145, 6).Value  = $ws.Cells.Item(147, 6).Value()   # F - Categoria ID
$ws.Cells.Item(145, 7).Value  = $ws.Cells.Item(147, 7).Value()   # G - Categoria
$ws.Cells.Item(145, 8).Value  = $ws.Cells.Item(147, 8).Value()   # H - Variedad
$ws.Cells.Item(145, 9).Value  = "Primera"                         # I - Calidad
$ws.Cells.Item(145, 10).Value = 133                                # J - Volumen
$ws.Cells.Item(145, 11).Value = 2000                               # K - Precio minimo
$ws.Cells.Item(145, 12).Value = 2300                               # L - Precio maximo
$ws.Cells.Item(145, 13).Value = 2151                               # M - Precio promedio ponderado
$ws.Cells.Item(145, 14).Value = $ws.Cells.Item(147, 14).Value()  # N - Unidad de comercializacion
$ws.Cells.Item(145, 15).Value = $ws.Cells.Item(147, 15).Value()  # O - Origen
$ws.Cells.Item(145, 16).Value = 717                                # P - Precio $/Kg
$ws.Cells.Item(145, 17).Value = $ws.Cells.Item(147, 17).Value()  # Q - Kg o Unidades
$ws.Cells.Item(145, 18).Value = $ws.Cells.Item(147, 18).Value()  # R - Clasificacion

$ws.Cells.Item(146, 1).Value  = $ws.Cells.Item(147, 1).Value()   # A - Mercado ID
$ws.Cells.Item(146, 2).Value  = $ws.Cells.Item(147, 2).Value()   # B - Mercado
$ws.Cells.Item(146, 3).Value  = $ws.Cells.Item(147, 3).Value()   # C - Region
$ws.Cells.Item(146, 4).Value  = 44452                             # D - Fecha
$ws.Cells.Item(146, 5).Value  = $ws.Cells.Item(147, 5).Value()   # E - Codreg
$ws.Cells.Item(146, 6).Value  = $ws.Cells.Item(147, 6).Value()   # F - Categoria ID
$ws.Cells.Item(146, 7).Value  = $ws.Cells.Item(147, 7).Value()   # G - Categoria
$ws.Cells.Item(146, 8).Value  = $ws.Cells.Item(147, 8).Value()   # H - Variedad
$ws.Cells.Item(146, 9).Value  = "Segunda"                         # I - Calidad
$ws.Cells.Item(146, 10).Value = 79                                 # J - Volumen
$ws.Cells.Item(146, 11).Value = 1500                               # K - Precio minimo
$ws.Cells.Item(146, 12).Value = 1800                               # L - Precio maximo
$ws.Cells.Item(146, 13).Value = 1648                               # M - Precio promedio ponderado
$ws.Cells.Item(146, 14).Value = $ws.Cells.Item(147, 14).Value()  # N - Unidad de comercializacion
$ws.Cells.Item(146, 15).Value = $ws.Cells.Item(147, 15).Value()  # O - Origen
$ws.Cells.Item(146, 16).Value = 549                                # P - Precio $/Kg
$ws.Cells.Item(146, 17).Value = $ws.Cells.Item(147, 17).Value()  # Q - Kg o Unidades
$ws.Cells.Item(146, 18).Value = $ws.Cells.Item(147, 18).Value()  # R - Clasificacion
